$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.988.17'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.88'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.85'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4249'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3656'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07227'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8414'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.53'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.821.41'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.650'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07063'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.274'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.48'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008742'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.174.80'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.124'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.80'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.059.71'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.975'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.82'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.229'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.16'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.51'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08709'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.174'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7342'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.902'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.089'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01939'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05207'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.204'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.867'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1682'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5106'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.511'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.48'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.954'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4725'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.73'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06318'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.647'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.59%  '
